$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2963.4
$ws.Range("I38").Value = 3227
$ws.Range("J38").Value = 1909
$ws.Range("K38").Value = 9681
$ws.Range("L38").Value = 5727
$ws.Range("M38").Value = -9309
$ws.Range("N38").Value = -6471

$ws.Range("H39").Value = 810
$ws.Range("I39").Value = 921.53845
$ws.Range("J39").Value = 520
$ws.Range("K39").Value = 2764.61535
$ws.Range("L39").Value = 1560
$ws.Range("M39").Value = -2468.61535
$ws.Range("N39").Value = -2152

$ws.Range("H42").Value = 1129.6364
$ws.Range("J42").Value = 179.6
$ws.Range("L42").Value = 538.8
$ws.Range("N42").Value = -998.8

$ws.Range("H43").Value = 585.7778
$ws.Range("I43").Value = 650
$ws.Range("J43").Value = 534.4
$ws.Range("K43").Value = 650
$ws.Range("L43").Value = 534.4
$ws.Range("M43").Value = -581
$ws.Range("N43").Value = -672.4

$ws.Range("H64").Value = 5264.091
$ws.Range("I64").Value = 3100
$ws.Range("J64").Value = 5480.5
$ws.Range("K64").Value = 3100
$ws.Range("L64").Value = 5480.5
$ws.Range("M64").Value = -2852
$ws.Range("N64").Value = -5976.5

$ws.Range("H67").Value = 5264.091
$ws.Range("I67").Value = 3100
$ws.Range("J67").Value = 5480.5
$ws.Range("K67").Value = 3100
$ws.Range("L67").Value = 5480.5
$ws.Range("M67").Value = -2242
$ws.Range("N67").Value = -7196.5

$ws.Range("H135").Value = 5588
$ws.Range("I135").Value = 6484.25
$ws.Range("K135").Value = 58358.25
$ws.Range("M135").Value = -55823.25

$ws.Range("H137").Value = 1356.5385
$ws.Range("I137").Value = 1000.5
$ws.Range("J137").Value = 1421.2727
$ws.Range("K137").Value = 3001.5
$ws.Range("L137").Value = 4263.8181
$ws.Range("M137").Value = -451.5
$ws.Range("N137").Value = -9363.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2343.1692
$ws.Range("I32").Value = 2251.3618
$ws.Range("J32").Value = 2582.889
$ws.Range("K32").Value = 2251.3618
$ws.Range("L32").Value = 2582.889
$ws.Range("M32").Value = -1964.3618
$ws.Range("N32").Value = -3156.889

$ws.Range("H45").Value = 2612.9333
$ws.Range("I45").Value = 1418
$ws.Range("J45").Value = 5002.8
$ws.Range("K45").Value = 1418
$ws.Range("L45").Value = 5002.8
$ws.Range("M45").Value = -1041
$ws.Range("N45").Value = -5756.8

$ws.Range("H61").Value = 2041
$ws.Range("I61").Value = 1641.75
$ws.Range("J61").Value = 3039.125
$ws.Range("K61").Value = 1641.75
$ws.Range("L61").Value = 3039.125
$ws.Range("M61").Value = -1429.75
$ws.Range("N61").Value = -3463.125

$ws.Range("H74").Value = 1187.0834
$ws.Range("I74").Value = 1278.3334
$ws.Range("J74").Value = 1095.8334
$ws.Range("K74").Value = 1278.3334
$ws.Range("L74").Value = 1095.8334
$ws.Range("M74").Value = -404.3334
$ws.Range("N74").Value = -2843.8334

$ws.Range("H77").Value = 1187.0834
$ws.Range("I77").Value = 1278.3334
$ws.Range("J77").Value = 1095.8334
$ws.Range("K77").Value = 6391.666999999999
$ws.Range("L77").Value = 5479.166999999999
$ws.Range("M77").Value = -2023.666999999999
$ws.Range("N77").Value = -14215.167

$ws.Range("H97").Value = 18518920
$ws.Range("I97").Value = 18518920
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 18518920
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -18518424
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 2168.7144
$ws.Range("I132").Value = 1564.0264
$ws.Range("K132").Value = 4692.0792
$ws.Range("M132").Value = -2162.0792

$ws.Range("H136").Value = 2041
$ws.Range("I136").Value = 1641.75
$ws.Range("J136").Value = 3039.125
$ws.Range("K136").Value = 4925.25
$ws.Range("L136").Value = 9117.375
$ws.Range("M136").Value = -2375.25
$ws.Range("N136").Value = -14217.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1844.7906
$ws.Range("I31").Value = 1103.037
$ws.Range("J31").Value = 3096.5
$ws.Range("K31").Value = 1103.037
$ws.Range("L31").Value = 3096.5
$ws.Range("M31").Value = -808.037
$ws.Range("N31").Value = -3686.5

$ws.Range("H34").Value = 1844.7906
$ws.Range("I34").Value = 1103.037
$ws.Range("J34").Value = 3096.5
$ws.Range("K34").Value = 1103.037
$ws.Range("L34").Value = 3096.5
$ws.Range("M34").Value = -901.037
$ws.Range("N34").Value = -3500.5

$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1317.9546
$ws.Range("I5").Value = 641.4286
$ws.Range("J5").Value = 1633.6666
$ws.Range("K5").Value = 1924.2858
$ws.Range("L5").Value = 4900.9998
$ws.Range("M5").Value = -1812.2858
$ws.Range("N5").Value = -5124.9998

$ws.Range("H68").Value = 2429.1482
$ws.Range("J68").Value = 1239.9354
$ws.Range("L68").Value = 3719.8062
$ws.Range("N68").Value = -5341.8062

$ws.Range("H71").Value = 2429.1482
$ws.Range("J71").Value = 1239.9354
$ws.Range("L71").Value = 11159.4186
$ws.Range("N71").Value = -19271.4186

$ws.Range("H135").Value = 1317.9546
$ws.Range("I135").Value = 641.4286
$ws.Range("J135").Value = 1633.6666
$ws.Range("K135").Value = 5772.8574
$ws.Range("L135").Value = 14702.9994
$ws.Range("M135").Value = -3237.8574
$ws.Range("N135").Value = -19772.9994

$ws.Range("H138").Value = 3507.8572
$ws.Range("I138").Value = 911
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 2733
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = 2407
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2462.1628
$ws.Range("I132").Value = 1985.5
$ws.Range("J132").Value = 4262.8887
$ws.Range("K132").Value = 5956.5
$ws.Range("L132").Value = 12788.6661
$ws.Range("M132").Value = -3426.5
$ws.Range("N132").Value = -17848.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 551.125
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 487
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 487
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -863

$ws.Range("H68").Value = 1701.625
$ws.Range("I68").Value = 1506.4546
$ws.Range("J68").Value = 2131
$ws.Range("K68").Value = 1506.4546
$ws.Range("L68").Value = 2131
$ws.Range("M68").Value = -757.4546
$ws.Range("N68").Value = -3629

$ws.Range("H71").Value = 1701.625
$ws.Range("I71").Value = 1506.4546
$ws.Range("J71").Value = 2131
$ws.Range("K71").Value = 7532.273
$ws.Range("L71").Value = 10655
$ws.Range("M71").Value = -3788.273
$ws.Range("N71").Value = -18143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1190.15
$ws.Range("I113").Value = 1396.8667
$ws.Range("J113").Value = 570
$ws.Range("K113").Value = 4190.6001
$ws.Range("L113").Value = 1710
$ws.Range("M113").Value = -2020.6001
$ws.Range("N113").Value = -6050

$ws.Range("H136").Value = 20897330
$ws.Range("I136").Value = 33433986
$ws.Range("J136").Value = 2903
$ws.Range("K136").Value = 100301958
$ws.Range("L136").Value = 8709
$ws.Range("M136").Value = -100299408
$ws.Range("N136").Value = -13809
